$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 210.66667
$ws.Range("I2").Value = 222
$ws.Range("J2").Value = 171
$ws.Range("K2").Value = 222
$ws.Range("L2").Value = 171
$ws.Range("M2").Value = -109
$ws.Range("N2").Value = -397
$ws.Range("H9").Value = 610.8823
$ws.Range("I9").Value = 697.61536
$ws.Range("J9").Value = 329
$ws.Range("K9").Value = 697.61536
$ws.Range("L9").Value = 329
$ws.Range("M9").Value = -528.61536
$ws.Range("N9").Value = -667
$ws.Range("H12").Value = 491.76923
$ws.Range("I12").Value = 528.4545000000001
$ws.Range("K12").Value = 528.4545000000001
$ws.Range("M12").Value = -358.4545000000001
$ws.Range("H19").Value = 1183.75
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 1183.75
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 1183.75
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -1533.75
$ws.Range("H33").Value = 324.86957
$ws.Range("J33").Value = 599
$ws.Range("L33").Value = 599
$ws.Range("N33").Value = -1057
$ws.Range("H40").Value = 2177.1428
$ws.Range("I40").Value = 1049.2
$ws.Range("J40").Value = 4997
$ws.Range("K40").Value = 1049.2
$ws.Range("L40").Value = 4997
$ws.Range("M40").Value = -874.2
$ws.Range("N40").Value = -5347
$ws.Range("H41").Value = 516.375
$ws.Range("I41").Value = 556.1667
$ws.Range("K41").Value = 556.1667
$ws.Range("M41").Value = -116.1667
$ws.Range("H70").Value = 3351.111
$ws.Range("I70").Value = 5009.2
$ws.Range("J70").Value = 1278.5
$ws.Range("K70").Value = 15027.6
$ws.Range("L70").Value = 3835.5
$ws.Range("M70").Value = -14757.6
$ws.Range("N70").Value = -4375.5
$ws.Range("H73").Value = 3351.111
$ws.Range("I73").Value = 5009.2
$ws.Range("J73").Value = 1278.5
$ws.Range("K73").Value = 15027.6
$ws.Range("L73").Value = 3835.5
$ws.Range("M73").Value = -14091.6
$ws.Range("N73").Value = -5707.5
$ws.Range("H74").Value = 5322
$ws.Range("I74").Value = 5322
$ws.Range("K74").Value = 5322
$ws.Range("M74").Value = -4386
$ws.Range("H76").Value = 3566
$ws.Range("I76").Value = 3899
$ws.Range("K76").Value = 3899
$ws.Range("M76").Value = -3584
$ws.Range("H77").Value = 5322
$ws.Range("I77").Value = 5322
$ws.Range("K77").Value = 26610
$ws.Range("M77").Value = -21930
$ws.Range("H79").Value = 3566
$ws.Range("I79").Value = 3899
$ws.Range("K79").Value = 3899
$ws.Range("M79").Value = -2807
$ws.Range("H137").Value = 1942.0834
$ws.Range("I137").Value = 1913.875
$ws.Range("J137").Value = 1998.5
$ws.Range("K137").Value = 5741.625
$ws.Range("L137").Value = 5995.5
$ws.Range("M137").Value = -3191.625
$ws.Range("N137").Value = -11095.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1394.7693
$ws.Range("I88").Value = 778.8333
$ws.Range("J88").Value = 1922.7142
$ws.Range("K88").Value = 778.8333
$ws.Range("L88").Value = 1922.7142
$ws.Range("M88").Value = -372.8333
$ws.Range("N88").Value = -2734.7142
$ws.Range("H91").Value = 1394.7693
$ws.Range("I91").Value = 778.8333
$ws.Range("J91").Value = 1922.7142
$ws.Range("K91").Value = 778.8333
$ws.Range("L91").Value = 1922.7142
$ws.Range("M91").Value = 625.1667
$ws.Range("N91").Value = -4730.7142
$ws.Range("H95").Value = 57879
$ws.Range("J95").Value = 57879
$ws.Range("L95").Value = 57879
$ws.Range("N95").Value = -63371
$ws.Range("H106").Value = 19300.334
$ws.Range("J106").Value = 19300.334
$ws.Range("L106").Value = 19300.334
$ws.Range("N106").Value = -21824.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2008.875
$ws.Range("J20").Value = 1287
$ws.Range("L20").Value = 1287
$ws.Range("N20").Value = -1781
$ws.Range("H26").Value = 11499.5
$ws.Range("I26").Value = 11499.5
$ws.Range("K26").Value = 11499.5
$ws.Range("M26").Value = -11207.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 50000
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("H62").Value = 2980.2222
$ws.Range("J62").Value = 2994
$ws.Range("L62").Value = 2994
$ws.Range("N62").Value = -4242
$ws.Range("H65").Value = 2980.2222
$ws.Range("J65").Value = 2994
$ws.Range("L65").Value = 14970
$ws.Range("N65").Value = -21210

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 3402.25
$ws.Range("J60").Value = 3870.2222
$ws.Range("L60").Value = 11610.6666
$ws.Range("N60").Value = -12112.6666
$ws.Range("H94").Value = 3006.6667
$ws.Range("J94").Value = 3027
$ws.Range("L94").Value = 9081
$ws.Range("N94").Value = -10433

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 321.16666
$ws.Range("I2").Value = 321.16666
$ws.Range("K2").Value = 321.16666
$ws.Range("M2").Value = -208.16666
$ws.Range("H70").Value = 7554.4614
$ws.Range("I70").Value = 6451.5
$ws.Range("K70").Value = 6451.5
$ws.Range("M70").Value = -6181.5
$ws.Range("H73").Value = 7554.4614
$ws.Range("I73").Value = 6451.5
$ws.Range("K73").Value = 6451.5
$ws.Range("M73").Value = -5515.5
$ws.Range("H80").Value = 2145.4614
$ws.Range("J80").Value = 2110.5715
$ws.Range("L80").Value = 2110.5715
$ws.Range("N80").Value = -4106.5715
$ws.Range("H83").Value = 2145.4614
$ws.Range("J83").Value = 2110.5715
$ws.Range("L83").Value = 10552.8575
$ws.Range("N83").Value = -20536.8575
$ws.Range("H92").Value = 13870.833
$ws.Range("J92").Value = 13870.833
$ws.Range("L92").Value = 13870.833
$ws.Range("N92").Value = -17614.833
$ws.Range("H101").Value = 25025
$ws.Range("J101").Value = 19550
$ws.Range("L101").Value = 19550
$ws.Range("N101").Value = -26040
$ws.Range("H102").Value = 3787.375
$ws.Range("I102").Value = 4114.2856
$ws.Range("J102").Value = 1499
$ws.Range("K102").Value = 4114.2856
$ws.Range("L102").Value = 1499
$ws.Range("M102").Value = -2492.2856
$ws.Range("N102").Value = -4743

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 8256
$ws.Range("I68").Value = 6999.5
$ws.Range("K68").Value = 6999.5
$ws.Range("M68").Value = -6250.5
$ws.Range("H71").Value = 8256
$ws.Range("I71").Value = 6999.5
$ws.Range("K71").Value = 34997.5
$ws.Range("M71").Value = -31253.5
$ws.Range("H82").Value = 2712
$ws.Range("I82").Value = 3232.6667
$ws.Range("J82").Value = 2399.6
$ws.Range("K82").Value = 3232.6667
$ws.Range("L82").Value = 2399.6
$ws.Range("M82").Value = -2871.6667
$ws.Range("N82").Value = -3121.6
$ws.Range("H85").Value = 2712
$ws.Range("I85").Value = 3232.6667
$ws.Range("J85").Value = 2399.6
$ws.Range("K85").Value = 3232.6667
$ws.Range("L85").Value = 2399.6
$ws.Range("M85").Value = -1984.6667
$ws.Range("N85").Value = -4895.6
$ws.Range("H97").Value = 39500
$ws.Range("J97").Value = 39500
$ws.Range("L97").Value = 39500
$ws.Range("N97").Value = -41482
$ws.Range("H101").Value = 12000
$ws.Range("J101").Value = 12000
$ws.Range("L101").Value = 12000
$ws.Range("N101").Value = -18490
$ws.Range("H106").Value = 42329.668
$ws.Range("J106").Value = 42329.668
$ws.Range("L106").Value = 42329.668
$ws.Range("N106").Value = -44853.668

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4798.9375
$ws.Range("I81").Value = 4178.5
$ws.Range("K81").Value = 8357
$ws.Range("M81").Value = -7296
$ws.Range("H84").Value = 4798.9375
$ws.Range("I84").Value = 4178.5
$ws.Range("K84").Value = 41785
$ws.Range("M84").Value = -36481
$ws.Range("H94").Value = 39450
$ws.Range("J94").Value = 39450
$ws.Range("L94").Value = 39450
$ws.Range("N94").Value = -41252
$ws.Range("H95").Value = 49750
$ws.Range("J95").Value = 49750
$ws.Range("L95").Value = 49750
$ws.Range("N95").Value = -55242
$ws.Range("H101").Value = 23499.5
$ws.Range("J101").Value = 23499.5
$ws.Range("L101").Value = 23499.5
$ws.Range("N101").Value = -29989.5
$ws.Range("H105").Value = 42650
$ws.Range("J105").Value = 42650
$ws.Range("L105").Value = 42650
$ws.Range("N105").Value = -49638
$ws.Range("H122").Value = 1499.8
$ws.Range("I122").Value = 1499.8
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4499.4
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2049.4
$ws.Range("N122").ClearContents()
